# Edit: insert two new weekly records (Primera/Segunda) for Coliflor at the top
# of the date-ordered block (rows 239-298), pushing the existing 60 rows down
# by two positions (to rows 241-300), and extending the sheet dimension to
# A1:R300.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 239 and 240; this shifts rows 239:298 down to 241:300
# and auto-extends the sheet dimension/UsedRange to A1:R300.
$ws.Range("A239:A240").EntireRow.Insert()

# Populate the new row 239 ("Primera")
$ws.Range("A239").Value = 5
$ws.Range("B239").Value = "Macroferia Regional de Talca"
$ws.Range("C239").Value = "Maule"
$ws.Range("D239").Value = 44785
$ws.Range("D239").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E239").Value = 7
$ws.Range("F239").Value = 100112008
$ws.Range("G239").Value = "Coliflor"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 3000
$ws.Range("K239").Value = 1000
$ws.Range("L239").Value = 1000
$ws.Range("M239").Value = 1000
$ws.Range("N239").Value = "`$/unidad"
$ws.Range("O239").Value = "Región del Maule"
$ws.Range("P239").Value = 1000
$ws.Range("Q239").Value = 1
$ws.Range("R239").Value = "Hortaliza"

# Populate the new row 240 ("Segunda")
$ws.Range("A240").Value = 5
$ws.Range("B240").Value = "Macroferia Regional de Talca"
$ws.Range("C240").Value = "Maule"
$ws.Range("D240").Value = 44785
$ws.Range("D240").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E240").Value = 7
$ws.Range("F240").Value = 100112008
$ws.Range("G240").Value = "Coliflor"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Segunda"
$ws.Range("J240").Value = 2000
$ws.Range("K240").Value = 800
$ws.Range("L240").Value = 800
$ws.Range("M240").Value = 800
$ws.Range("N240").Value = "`$/unidad"
$ws.Range("O240").Value = "Región del Maule"
$ws.Range("P240").Value = 800
$ws.Range("Q240").Value = 1
$ws.Range("R240").Value = "Hortaliza"
